# Applies the "clustersEmpleados" update:
#  1) Updates the K (Clusters) column values for a set of existing rows.
#  2) Appends a new employee row (row 272) with full data A272:K272.
#
# Matches the commit: "Funcionalidad y navegacion de editar nomina listo,
# ademas de la funcion de agregar empleado implementada."

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Update existing K-column (Clusters) values -------------------------
$kUpdates = @(
    @(2,0), @(3,2), @(7,2), @(9,2), @(15,4), @(16,4), @(17,0), @(18,2),
    @(22,2), @(24,2), @(30,4), @(31,4), @(32,0), @(33,2), @(37,2), @(39,2),
    @(45,4), @(46,4), @(47,0), @(48,2), @(52,2), @(54,2), @(60,4), @(61,4),
    @(62,0), @(63,2), @(67,2), @(69,2), @(75,4), @(76,4), @(77,0), @(78,2),
    @(82,2), @(84,2), @(90,4), @(91,4), @(92,0), @(93,2), @(97,2), @(99,2),
    @(105,4), @(106,4), @(107,0), @(108,2), @(112,2), @(114,2), @(120,4),
    @(121,4), @(122,0), @(123,2), @(127,2), @(129,2), @(135,4), @(136,4),
    @(137,0), @(138,2), @(142,2), @(144,2), @(150,4), @(151,4), @(152,0),
    @(153,2), @(157,2), @(159,2), @(165,4), @(166,4), @(167,0), @(168,2),
    @(172,2), @(174,2), @(180,4), @(181,4), @(182,0), @(183,2), @(187,2),
    @(189,2), @(195,4), @(196,4), @(197,0), @(198,2), @(202,2), @(204,2),
    @(210,4), @(211,4), @(212,0), @(213,2), @(217,2), @(219,2), @(225,4),
    @(226,4), @(227,0), @(228,2), @(232,2), @(234,2), @(240,4), @(241,4),
    @(242,0), @(243,2), @(247,2), @(249,2), @(255,4), @(256,4), @(257,0),
    @(258,2), @(262,2), @(264,2), @(270,4), @(271,4)
)

foreach ($pair in $kUpdates) {
    $row = $pair[0]
    $val = $pair[1]
    $ws.Cells.Item($row, 11).Value = $val
}

# --- 2) Add new employee row (row 272) --------------------------------------
$newRow = 272
$newValues = @(270, 10000, 3, 2, 12, 1, 1, 0, 1, 0, 1)
for ($col = 1; $col -le $newValues.Length; $col++) {
    $ws.Cells.Item($newRow, $col).Value = $newValues[$col - 1]
}

# Column A carries the bold/bordered/centered "id" cell format in every data
# row (same formatting as the header row). Copy it down from the row above
# so the new row matches the existing pattern exactly.
$ws.Cells.Item($newRow - 1, 1).Copy()
$ws.Cells.Item($newRow, 1).PasteSpecial(-4122)
